$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.44663052312135
$ws.Range("D2").Value = 3.744791562036846
$ws.Range("E2").Value = 8.987019149602336
$ws.Range("F2").Value = 51.19620335464731
$ws.Range("G2").Value = 3.775418369640712
$ws.Range("I2").Value = 34.51454239726278
$ws.Range("J2").Value = 9.114523741317687
$ws.Range("K2").Value = 18.38102186517671
$ws.Range("L2").Value = 13.47030010850402
$ws.Range("M2").Value = 17.88934705213561
$ws.Range("N2").Value = 24.90889892751992
$ws.Range("B3").Value = 15.39088126544116
$ws.Range("D3").Value = 3.71131348687254
$ws.Range("E3").Value = 8.940225338429739
$ws.Range("F3").Value = 51.20049217513596
$ws.Range("G3").Value = 3.778686819251128
$ws.Range("I3").Value = 34.57785010585479
$ws.Range("J3").Value = 9.101295316052243
$ws.Range("K3").Value = 18.24433696066139
$ws.Range("L3").Value = 13.46986750733029
$ws.Range("M3").Value = 17.89085616829107
$ws.Range("N3").Value = 24.9643250485044
$ws.Range("B4").Value = 15.35973310616355
$ws.Range("D4").Value = 3.690196349506404
$ws.Range("E4").Value = 8.910704185522027
$ws.Range("F4").Value = 51.21281854018044
$ws.Range("G4").Value = 3.780799692551164
$ws.Range("I4").Value = 34.6212957877345
$ws.Range("J4").Value = 9.092956631351747
$ws.Range("K4").Value = 18.16457686973193
$ws.Range("L4").Value = 13.47192399195532
$ws.Range("M4").Value = 17.89466763793109
$ws.Range("N4").Value = 25.00028460243926
$ws.Range("B5").Value = 15.34782515982143
$ws.Range("D5").Value = 3.681450852321535
$ws.Range("E5").Value = 8.898476294455172
$ws.Range("F5").Value = 51.22027855150262
$ws.Range("G5").Value = 3.781687461577929
$ws.Range("I5").Value = 34.64014986692067
$ws.Range("J5").Value = 9.089503971854478
$ws.Range("K5").Value = 18.13315171025139
$ws.Range("L5").Value = 13.47334670706726
$ws.Range("M5").Value = 17.89694771560725
$ws.Range("N5").Value = 25.01542388456189
$ws.Range("B6").Value = 15.34589555706213
$ws.Range("D6").Value = 3.679990240028092
$ws.Range("E6").Value = 8.896433939013331
$ws.Range("F6").Value = 51.22166442994286
$ws.Range("G6").Value = 3.781836493690254
$ws.Range("I6").Value = 34.64334998482413
$ws.Range("J6").Value = 9.088927359061957
$ws.Range("K6").Value = 18.12799952764359
$ws.Range("L6").Value = 13.47361827896364
$ws.Range("M6").Value = 17.89737025940982
$ws.Range("N6").Value = 25.01796709403786
$ws.Range("B7").Value = 15.35956931951026
$ws.Range("D7").Value = 3.690078970423923
$ws.Range("E7").Value = 8.910540075685027
$ws.Range("F7").Value = 51.21290928305001
$ws.Range("G7").Value = 3.780811556857965
$ws.Range("I7").Value = 34.62154540672337
$ws.Range("J7").Value = 9.092910288773382
$ws.Range("K7").Value = 18.16414865672486
$ws.Range("L7").Value = 13.47194081101674
$ws.Range("M7").Value = 17.89469544307054
$ws.Range("N7").Value = 25.00048680948982
$ws.Range("B8").Value = 15.42677484288486
$ws.Range("D8").Value = 3.733365501107094
$ws.Range("E8").Value = 8.971048968200765
$ws.Range("F8").Value = 51.19567037040268
$ws.Range("G8").Value = 3.776523380422404
$ws.Range("I8").Value = 34.53542106645324
$ws.Range("J8").Value = 9.110007386952148
$ws.Range("K8").Value = 18.33304730244
$ws.Range("L8").Value = 13.46966952224831
$ws.Range("M8").Value = 17.88926958282876
$ws.Range("N8").Value = 24.92761020606874
$ws.Range("B9").Value = 15.58253485739039
$ws.Range("D9").Value = 3.813749947244634
$ws.Range("E9").Value = 9.08341141238833
$ws.Range("F9").Value = 51.23875835299073
$ws.Range("G9").Value = 3.768951339895304
$ws.Range("I9").Value = 34.40285787288181
$ws.Range("J9").Value = 9.141826280102947
$ws.Range("K9").Value = 18.69587775238925
$ws.Range("L9").Value = 13.48359851745319
$ws.Range("M9").Value = 17.90144732642597
$ws.Range("N9").Value = 24.79996081446794
$ws.Range("B10").Value = 15.71088136021323
$ws.Range("D10").Value = 3.870008273022046
$ws.Range("E10").Value = 9.162089970379119
$ws.Range("F10").Value = 51.31721470696336
$ws.Range("G10").Value = 3.763892505459354
$ws.Range("I10").Value = 34.32764892224134
$ws.Range("J10").Value = 9.164174925395807
$ws.Range("K10").Value = 18.97970036868769
$ws.Range("L10").Value = 13.50497079501259
$ws.Range("M10").Value = 17.92419496735716
$ws.Range("N10").Value = 24.71543267399244
$ws.Range("B11").Value = 15.77212402193122
$ws.Range("D11").Value = 3.894982109039008
$ws.Range("E11").Value = 9.197034092249618
$ws.Range("F11").Value = 51.36302938806975
$ws.Range("G11").Value = 3.761699356259109
$ws.Range("I11").Value = 34.29825895420102
$ws.Range("J11").Value = 9.174121303697124
$ws.Range("K11").Value = 19.11210524256136
$ws.Range("L11").Value = 13.51709272731907
$ws.Range("M11").Value = 17.93751130772368
$ws.Range("N11").Value = 24.67897850917602
$ws.Range("B12").Value = 15.79571080382149
$ws.Range("D12").Value = 3.904349108274411
$ws.Range("E12").Value = 9.210143930605494
$ws.Range("F12").Value = 51.38182859728843
$ws.Range("G12").Value = 3.760884321243998
$ws.Range("I12").Value = 34.28782371485751
$ws.Range("J12").Value = 9.177856266362197
$ws.Range("K12").Value = 19.16267693525893
$ws.Range("L12").Value = 13.52202591889299
$ws.Range("M12").Value = 17.94297774431601
$ws.Range("N12").Value = 24.66546090587475
$ws.Range("B13").Value = 15.79061362859696
$ws.Range("D13").Value = 3.902335780135876
$ws.Range("E13").Value = 9.207325966025104
$ws.Range("F13").Value = 51.37771545739636
$ws.Range("G13").Value = 3.761059167282785
$ws.Range("I13").Value = 34.29004025181878
$ws.Range("J13").Value = 9.177053274762052
$ws.Range("K13").Value = 19.15176681930039
$ws.Range("L13").Value = 13.52094825559247
$ws.Range("M13").Value = 17.94178165053303
$ws.Range("N13").Value = 24.66835941639043
$ws.Range("B14").Value = 15.77405668406579
$ws.Range("D14").Value = 3.895754548828633
$ws.Range("E14").Value = 9.198115111754639
$ws.Range("F14").Value = 51.36454698340368
$ws.Range("G14").Value = 3.76163199343539
$ws.Range("I14").Value = 34.29738652698077
$ws.Range("J14").Value = 9.174429212136049
$ws.Range("K14").Value = 19.11625736637258
$ws.Range("L14").Value = 13.51749172554207
$ws.Range("M14").Value = 17.93795255888445
$ws.Range("N14").Value = 24.67786066292108
$ws.Range("B15").Value = 15.7639661191248
$ws.Range("D15").Value = 3.891711594295121
$ws.Range("E15").Value = 9.19245718554353
$ws.Range("F15").Value = 51.35666957990902
$ws.Range("G15").Value = 3.761984877146554
$ws.Range("I15").Value = 34.30197674072923
$ws.Range("J15").Value = 9.172817797946456
$ws.Range("K15").Value = 19.09456195288117
$ws.Range("L15").Value = 13.5154190850185
$ws.Range("M15").Value = 17.93566223027133
$ws.Range("N15").Value = 24.68371778495001
$ws.Range("B16").Value = 15.70693520490936
$ws.Range("D16").Value = 3.86836363383443
$ws.Range("E16").Value = 9.159789164069348
$ws.Range("F16").Value = 51.31442387992114
$ws.Range("G16").Value = 3.764038001485914
$ws.Range("I16").Value = 34.32966672666207
$ws.Range("J16").Value = 9.163520488564652
$ws.Range("K16").Value = 18.97111019643027
$ws.Range("L16").Value = 13.50422670360725
$ws.Range("M16").Value = 17.9233841922849
$ws.Range("N16").Value = 24.71785520089637
$ws.Range("B17").Value = 15.67266966666542
$ws.Range("D17").Value = 3.853881226584127
$ws.Range("E17").Value = 9.139530758494777
$ws.Range("F17").Value = 51.29109695508943
$ws.Range("G17").Value = 3.765325161640906
$ws.Range("I17").Value = 34.34788934444882
$ws.Range("J17").Value = 9.1577606094821
$ws.Range("K17").Value = 18.89619116775456
$ws.Range("L17").Value = 13.49797366195706
$ws.Range("M17").Value = 17.91661021808315
$ws.Range("N17").Value = 24.73930878190276
$ws.Range("B18").Value = 15.65323092862458
$ws.Range("D18").Value = 3.845493192192642
$ws.Range("E18").Value = 9.127798976922286
$ws.Range("F18").Value = 51.27863330743632
$ws.Range("G18").Value = 3.766075685956421
$ws.Range("I18").Value = 34.35882444392598
$ws.Range("J18").Value = 9.154426951251125
$ws.Range("K18").Value = 18.85341211739819
$ws.Range("L18").Value = 13.49460308400071
$ws.Range("M18").Value = 17.9129935984142
$ws.Range("N18").Value = 24.75183641628313
$ws.Range("B19").Value = 15.64669611748618
$ws.Range("D19").Value = 3.842643197749148
$ws.Range("E19").Value = 9.12381313989839
$ws.Range("F19").Value = 51.27457722907967
$ws.Range("G19").Value = 3.766331552348738
$ws.Range("I19").Value = 34.36260482508951
$ws.Range("J19").Value = 9.153294667484074
$ws.Range("K19").Value = 18.83898274596534
$ws.Range("L19").Value = 13.49350074309081
$ws.Range("M19").Value = 17.9118171803546
$ws.Range("N19").Value = 24.75611038110421
$ws.Range("B20").Value = 15.6762894679924
$ws.Range("D20").Value = 3.855428928721789
$ws.Range("E20").Value = 9.141695553816813
$ws.Range("F20").Value = 51.29348151242426
$ws.Range("G20").Value = 3.765187087894045
$ws.Range("I20").Value = 34.34590253122779
$ws.Range("J20").Value = 9.158375902682874
$ws.Range("K20").Value = 18.90413438690345
$ws.Range("L20").Value = 13.49861593313672
$ws.Range("M20").Value = 17.91730240357719
$ws.Range("N20").Value = 24.7370055477173
$ws.Range("B21").Value = 15.77890925122785
$ws.Range("D21").Value = 3.897690068818967
$ws.Range("E21").Value = 9.200823901464107
$ws.Range("F21").Value = 51.3683755772205
$ws.Range("G21").Value = 3.761463321541375
$ws.Range("I21").Value = 34.29520990603437
$ws.Range("J21").Value = 9.175200817141958
$ws.Range("K21").Value = 19.12667593324296
$ws.Range("L21").Value = 13.51849770519748
$ws.Range("M21").Value = 17.93906577878038
$ws.Range("N21").Value = 24.6750621370026
$ws.Range("B22").Value = 15.84827503419286
$ws.Range("D22").Value = 3.924784779313053
$ws.Range("E22").Value = 9.238751833229527
$ws.Range("F22").Value = 51.42577247297741
$ws.Range("G22").Value = 3.759119714922838
$ws.Range("I22").Value = 34.26612531163347
$ws.Range("J22").Value = 9.186013241434047
$ws.Range("K22").Value = 19.27462337684365
$ws.Range("L22").Value = 13.53348896195412
$ws.Range("M22").Value = 17.95575819694798
$ws.Range("N22").Value = 24.63624998284153
$ws.Range("B23").Value = 15.81104818625914
$ws.Range("D23").Value = 3.910372261789951
$ws.Range("E23").Value = 9.218574787804403
$ws.Range("F23").Value = 51.39436775098124
$ws.Range("G23").Value = 3.76036232738636
$ws.Range("I23").Value = 34.28127793683309
$ws.Range("J23").Value = 9.180259203637101
$ws.Range("K23").Value = 19.19544553106051
$ws.Range("L23").Value = 13.5253058723646
$ws.Range("M23").Value = 17.94662431284205
$ws.Range("N23").Value = 24.65681199208718
$ws.Range("B24").Value = 15.67465214174906
$ws.Range("D24").Value = 3.854729405349505
$ws.Range("E24").Value = 9.140717114429464
$ws.Range("F24").Value = 51.29240050282598
$ws.Range("G24").Value = 3.765249478263039
$ws.Range("I24").Value = 34.34679934078505
$ws.Range("J24").Value = 9.158097797881046
$ws.Range("K24").Value = 18.90054234210225
$ws.Range("L24").Value = 13.49832486325812
$ws.Range("M24").Value = 17.91698860088783
$ws.Range("N24").Value = 24.73804623651983
$ws.Range("B25").Value = 15.5379009064711
$ws.Range("D25").Value = 3.792491516486225
$ws.Range("E25").Value = 9.053692076551567
$ws.Range("F25").Value = 51.21887331578433
$ws.Range("G25").Value = 3.770910781560374
$ws.Range("I25").Value = 34.43482633466157
$ws.Range("J25").Value = 9.133399992370892
$ws.Range("K25").Value = 18.59454064151442
$ws.Range("L25").Value = 13.47786836184341
$ws.Range("M25").Value = 17.8957203475713
$ws.Range("N25").Value = 24.83286453946465
